$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose content participates in the shuffle
$rows = @(15,16,17,18,19,20,21,22,25,26,27,28,29,30,34,35,43,44)
$cols = @("A","B","D","E","F","G","H","M","Q","R")

# Snapshot current ("before") values for every moving cell
$data = @{}
foreach ($r in $rows) {
  foreach ($c in $cols) {
    $addr = "$c$r"
    $data[$addr] = $ws.Range($addr).Value()
  }
}

# Destination row -> source row (content provenance), per the target diff
$map = @{}
$map[15] = 16
$map[16] = 17
$map[17] = 18
$map[18] = 19
$map[19] = 15
$map[20] = 22
$map[21] = 20
$map[22] = 21
$map[25] = 27
$map[26] = 28
$map[27] = 25
$map[28] = 26
$map[29] = 30
$map[30] = 29
$map[34] = 35
$map[35] = 34
$map[43] = 44
$map[44] = 43

# Write snapshotted source-row content into each destination row
foreach ($destRow in $rows) {
  $srcRow = $map[$destRow]
  foreach ($c in $cols) {
    $srcAddr = "$c$srcRow"
    $destAddr = "$c$destRow"
    $val = $data[$srcAddr]
    if ($val -eq $null) {
      $ws.Range($destAddr).Value = ""
    } else {
      $ws.Range($destAddr).Value = $val
    }
  }
}

Write-Output "done"